$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: new participant entry "xiaoyi_20251202_134616"
# Force text storage for the SmartScore "numeric-looking" columns so that
# values like "0.520" keep their trailing zero (stored as text, not a number).
foreach ($addr in @("I29", "L29", "O29", "R29", "U29", "X29", "AA29", "AD29", "AG29")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A29").Value = "xiaoyi_20251202_134616"
# B29: left blank (source row has an empty text cell here; the empty-string
#   COM Value assignment deletes the cell outright in this runtime, so it is
#   simply not written -- net effect is the same blank cell).
$ws.Range("C29").Value = "xiaoyi"
$ws.Range("D29").Value = 26
$ws.Range("E29").Value = "Female"
$ws.Range("F29").Value = "2025-12-02 13:46:17"
$ws.Range("G29").Value = "{
  `"portion`": 0.4,
  `"diet`": 0.7142857142857143,
  `"salt`": 0.2,
  `"fat`": 0.8,
  `"natural`": 0.2,
  `"convenience`": 1.0,
  `"price`": 1.0
}"
$ws.Range("H29").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("I29").Value = "0.591"
$ws.Range("J29").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("K29").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("L29").Value = "0.532"
$ws.Range("M29").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("N29").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O29").Value = "0.518"
$ws.Range("P29").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("Q29").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("R29").Value = "0.607"
$ws.Range("S29").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"
$ws.Range("T29").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U29").Value = "0.520"
$ws.Range("V29").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("W29").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("X29").Value = "0.443"
$ws.Range("Y29").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("Z29").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AA29").Value = "0.718"
$ws.Range("AB29").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AC29").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AD29").Value = "0.705"
$ws.Range("AE29").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"
$ws.Range("AF29").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AG29").Value = "0.673"
$ws.Range("AH29").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
